$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated worker / period data table (Data/EC/NIT-9009959584.xlsx)
# Each worker now appears with period 1801 before period 1712, and the
# three workers are reordered (Marilyn, Martha, Rosa). Rosa's mora/salario
# values also changed from 1023000 to 2630000.
$data = @(
    @{ Row = 16; Doc = "1143327174"; Nombre = "MARILYN CANTILLO AVILA";               Periodo = "1801"; Mora = 29509; Salario = 737717  },
    @{ Row = 17; Doc = "1143327174"; Nombre = "MARILYN CANTILLO AVILA";               Periodo = "1712"; Mora = 29509; Salario = 737717  },
    @{ Row = 18; Doc = "1148434693"; Nombre = "MARTHA IRINA TEJEDOR PANZA";           Periodo = "1801"; Mora = 29509; Salario = 737717  },
    @{ Row = 19; Doc = "1148434693"; Nombre = "MARTHA IRINA TEJEDOR PANZA";           Periodo = "1712"; Mora = 29509; Salario = 737717  },
    @{ Row = 20; Doc = "1128054808"; Nombre = "ROSA ALEJANDRA CARRASQUILLA RONCALLO"; Periodo = "1801"; Mora = 40920; Salario = 2630000 },
    @{ Row = 21; Doc = "1128054808"; Nombre = "ROSA ALEJANDRA CARRASQUILLA RONCALLO"; Periodo = "1712"; Mora = 40920; Salario = 2630000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("C$r").Value = $item.Doc
    $ws.Range("D$r").Value = $item.Nombre
    $ws.Range("E$r").Value = $item.Periodo
    $ws.Range("F$r").Value = $item.Mora
    $ws.Range("G$r").Value = $item.Salario
}

# Best-effort re-fit of the "bestFit" columns (B:J) to reflect the new
# (slightly wider) content, matching the recalculated widths Excel stored.
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
